# "explode sku unit on import"
# Add a new "value" column (G) to the sample import sheet and update the
# sample row's subcategory/category/sku/type values, plus give it a sample
# numeric "value" so the import demonstrates exploding a multi-sku cell
# (e.g. "skus1,skus2") into per-sku value rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header for column G. Set this before touching row 2 so the shared
# string table keeps the same ordering as the authored workbook.
$ws.Range("E1").Value = "sku"
$ws.Range("F1").Value = "type"
$ws.Range("G1").Value = "value"

# Updated sample data row.
$ws.Range("C2").Value = "sub1"
$ws.Range("D2").Value = "cat1"
$ws.Range("E2").Value = "skus1,skus2"
$ws.Range("F2").Value = "fast moving"
$ws.Range("G2").Value = 55

# Leave the selection on the newly added cell, matching where the author's
# cursor ended up after adding the column.
$ws.Range("G2").Select()
